# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 357
$ws1.Range("F6").Value = 5321
$ws1.Range("F8").Value = 5365
$ws1.Range("F9").Value = 624
$ws1.Range("F11").Value = 1371

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 357
$ws4.Range("F7").Value = 5321
$ws4.Range("F9").Value = 5365
$ws4.Range("F10").Value = 624
$ws4.Range("F12").Value = 1371
